# Corona data: remove "Confirmados_real" helper column (F), and add a new
# daily data row (row 27) below the existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column F (formulas + header "Confirmados_real") entirely.
$ws.Range("F1:F26").EntireColumn.Delete() | Out-Null

# Add the new data row (row 27).
$ws.Range("A27").Value = 86
$ws.Range("B27").Value = 22257
$ws.Range("C27").Value = 3544
$ws.Range("D27").Value = 43
$ws.Range("E27").Value = 60

# Update selection to match target state.
$ws.Range("J19").Select() | Out-Null
